$wb = $excel.ActiveWorkbook

# --- Sheet 1: "PI hours" ---
# Columns are currently: A=index, B=name, C=hours, D=dept
# New layout: A=index, B=name, C=hours, D=percentage, E=dept
$ws1 = $wb.Worksheets.Item("PI hours")

# Shift the "dept" column (D) one column to the right, into E, before
# inserting the new "percentage" column into D.
$ws1.Range("D1:D3").Cut($ws1.Range("E1:E3"))

# Header for the new column (copy the header style/border from column C)
$ws1.Range("D1").Value = "percentage"
$ws1.Range("C1").Copy()
$ws1.Range("D1").PasteSpecial(-4122)

# Percentage values = hours / total hours (this sheet) * 100
$ws1.Range("D2").Value = 60
$ws1.Range("D3").Value = 40

# --- Sheet 2: "dept hours" ---
# Columns are currently: A=index, B=dept, C=hours
# New layout: A=index, B=dept, C=hours, D=percentage
$ws2 = $wb.Worksheets.Item("dept hours")

$ws2.Range("D1").Value = "percentage"
$ws2.Range("C1").Copy()
$ws2.Range("D1").PasteSpecial(-4122)

# Percentage values = hours / total hours (this sheet) * 100
$ws2.Range("D2").Value = 41.66666666666666
$ws2.Range("D3").Value = 25
$ws2.Range("D4").Value = 16.66666666666667
$ws2.Range("D5").Value = 16.66666666666667
